# Bank of America test workbook: add two new tabs (WealthManagement,
# SmallBusiness) of Assertions test data next to the existing ContactUs tab.

$wb = $excel.ActiveWorkbook

$contactUs = $wb.Worksheets.Item("ContactUs")

# --- WealthManagement: inserted right after ContactUs -------------------
$wealth = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $contactUs)
$wealth.Name = "WealthManagement"
$wealth.Range("A1").Value = "Assertions"
$wealth.Range("A2").Value = "Thank you. A representative will be in contact shortly."
$wealth.Columns.Item(1).ColumnWidth = 46.08984375

# --- SmallBusiness: inserted right after WealthManagement ---------------
$smallBiz = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wealth)
$smallBiz.Name = "SmallBusiness"
$smallBiz.Range("A1").Value = "Assertions"
$smallBiz.Range("A2").Value = "Cash flow"
$smallBiz.Range("A3").Value = "Credit and funding"
$smallBiz.Range("A4").Value = "Taxes"
$smallBiz.Range("A5").Value = "Retirement"
$smallBiz.Range("A6").Value = "HR"
$smallBiz.Range("A7").Value = "Customer relations"
$smallBiz.Range("A8").Value = "Business strategy"
$smallBiz.Range("A9").Value = "News"
$smallBiz.Range("A10").Value = "Industry trends"
$smallBiz.Range("A11").Value = "Women Entrepreneurs"
$smallBiz.Range("A12").Value = "Small Business Spotlight"
$smallBiz.Range("A13").Value = "Heartbeat of Main Street"

# A handful of the topic rows wrap onto multiple lines in the sheet.
$smallBiz.Range("A8").WrapText = $true
$smallBiz.Range("A10").WrapText = $true
$smallBiz.Range("A12").WrapText = $true

$smallBiz.Columns.Item(1).ColumnWidth = 21.7265625

# --- View state: restore each sheet's last-used selection/active cell ---
$null = $wealth.Range("B28").Select()
$null = $smallBiz.Range("A14").Select()

# SmallBusiness ends up the active (selected) tab.
$null = $smallBiz.Activate()
